$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95-97 down to 96-98
$ws.Rows("95:95").Insert()

# Populate the newly inserted row 95 with its data
$ws.Range("A95").Value = 10
$ws.Range("B95").Value = "Vega Modelo de Temuco"
$ws.Range("C95").Value = "La Araucanía"
$ws.Range("D95").Value = 44509
$ws.Range("E95").Value = 9
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100107
$ws.Range("H95").Value = "Otros"
$ws.Range("I95").Value = 100107002
$ws.Range("J95").Value = "Chirimoya"
$ws.Range("K95").Value = "Cultivar IV Región"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 25
$ws.Range("N95").Value = 3000
$ws.Range("O95").Value = 3000
$ws.Range("P95").Value = 3000
$ws.Range("Q95").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R95").Value = "Provincia del Elquí"
$ws.Range("S95").Value = 3000
$ws.Range("T95").Value = 1
